$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# Remove the duplicated "Contact" row (row 11 duplicates row 10) so the
# sheet shrinks from 21 to 20 rows, same as the rest of the rows shifting up.
$ws1.Rows.Item(11).Delete()

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# New publication date
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher is now populated
$ws1.Range("B9").Value = "Alvearie Team"

# The old "Contact" row becomes a "Jurisdiction" row
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item(2)

# Root extension row: give it the real Short/Definition text instead of the
# generic Extension placeholder text.
$ws2.Range("K2").Value = "Parent Organization Hierarchy Level Code"
$ws2.Range("L2").Value = "Numeric level of the parent within the organinzational hierarchy"

# The "Short" column (K) now needs to be a bit wider to fit the longer text.
$ws2.Columns.Item(11).ColumnWidth = 38.8
